$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; existing rows 19-22 shift down to 20-23
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 45194
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100107
$ws.Cells.Item(19, 8).Value = "Otros"
$ws.Cells.Item(19, 9).Value = 100107002
$ws.Cells.Item(19, 10).Value = "Chirimoya"
$ws.Cells.Item(19, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 80
$ws.Cells.Item(19, 14).Value = 22000
$ws.Cells.Item(19, 15).Value = 22000
$ws.Cells.Item(19, 16).Value = 22000
$ws.Cells.Item(19, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 19).Value = 2200
$ws.Cells.Item(19, 20).Value = 10
